$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new product row (row 8) that was previously empty (only E8 had formatting).
$ws.Range("A8").Value = "Bel Paça Lastik  Pantolon"
$ws.Range("B8").Value = "350 Tl"
$ws.Range("C8").Value = "Jeans"
$ws.Range("D8").Value = "ANTRASİT.jpg"
$ws.Range("E8").Value = "%98 pamuk içeriği ile nefes alabilen yapıda, cildinize nazik dokunuşlar sunar.Bağcıklı kapama şekliyle kişisel zevke göre ayarlama imkanı taşır.Lastikli bel detayı sayesinde  rahatlığından ödün vermez.34-46 Beden aralığı mevcuttur."

# Update the active view/selection as it was left by the author (scrolled right, cell E8 selected).
[void]$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollColumn = 2
